# Apply the StructureDefinition-nursing-home-indicator.xlsx update:
#  - Rebrand from "ibm.com/Alvearie Team" to "linuxforhealth.org/LinuxForHealth Team"
#  - Bump Version 7.0.0 -> 8.0.0 and Date to the new publish timestamp
#  - Clear the stray duplicated ele-1/ext-1 constraint text that had been
#    sitting on the top-level Extension row (it correctly lives on the
#    Extension.extension row already).

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/nursing-home-indicator"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 is the top-level "Extension" element; column AI is "Constraint(s)".
# It incorrectly duplicated the ele-1/ext-1 invariant text that belongs to
# the "Extension.extension" row (row 4) only - clear it here.
$elements.Range("AI2").Value = ""
